$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("E4").Value = 21

# Row 17
$ws.Range("E17").Value = 97

# Row 25
$ws.Range("E25").Value = 16

# Row 26
$ws.Range("E26").Value = 25

# Row 32
$ws.Range("E32").Value = 16

# Row 33
$ws.Range("E33").Value = 36

# Row 36
$ws.Range("E36").Value = 89
$ws.Range("F36").Value = 35
$ws.Range("H36").Value = 35

# Row 42
$ws.Range("E42").Value = 28
$ws.Range("F42").Value = 11
$ws.Range("H42").Value = 11

# Row 48
$ws.Range("E48").Value = 25
$ws.Range("F48").Value = 13
$ws.Range("H48").Value = 13

# Row 49
$ws.Range("E49").Value = 55
$ws.Range("F49").Value = 31
$ws.Range("H49").Value = 31

# Row 55
$ws.Range("E55").Value = 7

# Row 62
$ws.Range("E62").Value = 38

# Row 64
$ws.Range("E64").Value = 30

# Row 70
$ws.Range("E70").Value = 37

# Row 82
$ws.Range("E82").Value = 14

# Row 85
$ws.Range("E85").Value = 6
$ws.Range("F85").Value = 4
$ws.Range("H85").Value = 4

# Row 88
$ws.Range("E88").Value = 19
